# Auto-generated edit script: updates cached leve-profit values in Sheets/Titan_Profits.xlsx
# Source data: scheduled runner refresh of currentAveragePrice / profit columns (H-N) per leve row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 104.9
$ws.Range("I33").Value = 99.611115
$ws.Range("K33").Value = 99.611115
$ws.Range("M33").Value = 129.388885
$ws.Range("H113").Value = 1166.25
$ws.Range("I113").Value = 1184.7222
$ws.Range("K113").Value = 1184.7222
$ws.Range("M113").Value = 2069.2778
$ws.Range("H137").Value = 45455860
$ws.Range("I137").Value = 62500830
$ws.Range("J137").Value = 2598.6667
$ws.Range("K137").Value = 187502490
$ws.Range("L137").Value = 7796.000100000001
$ws.Range("M137").Value = -187499940
$ws.Range("N137").Value = -12896.0001
$ws.Range("H138").Value = 4119578
$ws.Range("I138").Value = 1468844.5
$ws.Range("J138").Value = 5466672
$ws.Range("K138").Value = 4406533.5
$ws.Range("L138").Value = 16400016
$ws.Range("M138").Value = -4401393.5
$ws.Range("N138").Value = -16410296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17260.07
$ws.Range("I32").Value = 3057.6562
$ws.Range("K32").Value = 3057.6562
$ws.Range("M32").Value = -2770.6562
$ws.Range("H74").Value = 6342
$ws.Range("I74").Value = 1170.1177
$ws.Range("K74").Value = 1170.1177
$ws.Range("M74").Value = -296.1177
$ws.Range("H77").Value = 6342
$ws.Range("I77").Value = 1170.1177
$ws.Range("K77").Value = 5850.5885
$ws.Range("M77").Value = -1482.5885
$ws.Range("H102").Value = 3247.8948
$ws.Range("I102").Value = 3593.1538
$ws.Range("J102").Value = 2499.8333
$ws.Range("K102").Value = 3593.1538
$ws.Range("L102").Value = 2499.8333
$ws.Range("M102").Value = -1971.1538
$ws.Range("N102").Value = -5743.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2355.3076
$ws.Range("I94").Value = 2451.125
$ws.Range("J94").Value = 2202
$ws.Range("K94").Value = 2451.125
$ws.Range("L94").Value = 2202
$ws.Range("M94").Value = -2000.125
$ws.Range("N94").Value = -3104
$ws.Range("H105").Value = 248911.48
$ws.Range("I105").Value = 5718.207
$ws.Range("J105").Value = 836628.5600000001
$ws.Range("K105").Value = 5718.207
$ws.Range("L105").Value = 836628.5600000001
$ws.Range("M105").Value = -3971.207
$ws.Range("N105").Value = -840122.5600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 64249.875
$ws.Range("I16").Value = 101240
$ws.Range("J16").Value = 2599.6667
$ws.Range("K16").Value = 101240
$ws.Range("L16").Value = 2599.6667
$ws.Range("M16").Value = -100953
$ws.Range("N16").Value = -3173.6667
$ws.Range("H86").Value = 17245142
$ws.Range("I86").Value = 22731496
$ws.Range("K86").Value = 22731496
$ws.Range("M86").Value = -22730373
$ws.Range("H89").Value = 17245142
$ws.Range("I89").Value = 22731496
$ws.Range("K89").Value = 113657480
$ws.Range("M89").Value = -113651864
$ws.Range("H113").Value = 64249.875
$ws.Range("I113").Value = 101240
$ws.Range("J113").Value = 2599.6667
$ws.Range("K113").Value = 101240
$ws.Range("L113").Value = 2599.6667
$ws.Range("M113").Value = -99070
$ws.Range("N113").Value = -6939.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2933
$ws.Range("I86").Value = 3710.6667
$ws.Range("K86").Value = 11132.0001
$ws.Range("M86").Value = -9946.000100000001
$ws.Range("H89").Value = 2933
$ws.Range("I89").Value = 3710.6667
$ws.Range("K89").Value = 33396.0003
$ws.Range("M89").Value = -27468.0003
$ws.Range("H113").Value = 10000524
$ws.Range("J113").Value = 13889349
$ws.Range("L113").Value = 41668047
$ws.Range("N113").Value = -41672387

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 12334.667
$ws.Range("I4").Value = 25004
$ws.Range("J4").Value = 6000
$ws.Range("K4").Value = 25004
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = -24891
$ws.Range("N4").Value = -6226
$ws.Range("H28").Value = 12334.667
$ws.Range("I28").Value = 25004
$ws.Range("J28").Value = 6000
$ws.Range("K28").Value = 25004
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = -24772
$ws.Range("N28").Value = -6464
$ws.Range("H29").Value = 16007.2
$ws.Range("J29").Value = 7509
$ws.Range("L29").Value = 7509
$ws.Range("N29").Value = -8099
$ws.Range("H37").Value = 12334.667
$ws.Range("I37").Value = 25004
$ws.Range("J37").Value = 6000
$ws.Range("K37").Value = 25004
$ws.Range("L37").Value = 6000
$ws.Range("M37").Value = -24897
$ws.Range("N37").Value = -6214
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 1485.7142
$ws.Range("J68").Value = 3200
$ws.Range("K68").Value = 1485.7142
$ws.Range("L68").Value = 3200
$ws.Range("M68").Value = -736.7141999999999
$ws.Range("N68").Value = -4698
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 1485.7142
$ws.Range("J71").Value = 3200
$ws.Range("K71").Value = 7428.571
$ws.Range("L71").Value = 16000
$ws.Range("M71").Value = -3684.571
$ws.Range("N71").Value = -23488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 18666.666
$ws.Range("J9").Value = 3000
$ws.Range("L9").Value = 3000
$ws.Range("N9").Value = -3280
$ws.Range("H11").Value = 18666.334
$ws.Range("J11").Value = 2999.5
$ws.Range("L11").Value = 2999.5
$ws.Range("N11").Value = -3283.5
$ws.Range("H12").Value = 26500
$ws.Range("J12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("N12").Value = -3284
$ws.Range("H14").Value = 12514000
$ws.Range("I14").Value = 50000
$ws.Range("J14").Value = 16668667
$ws.Range("K14").Value = 50000
$ws.Range("L14").Value = 16668667
$ws.Range("M14").Value = -49832
$ws.Range("N14").Value = -16669003
$ws.Range("H17").Value = 4260.5557
$ws.Range("I17").Value = 4620.7144
$ws.Range("K17").Value = 4620.7144
$ws.Range("M17").Value = -4448.7144
$ws.Range("H24").Value = 50001500
$ws.Range("J24").Value = 3000
$ws.Range("L24").Value = 3000
$ws.Range("N24").Value = -3460
$ws.Range("H31").Value = 5750
$ws.Range("J31").Value = 5750
$ws.Range("L31").Value = 5750
$ws.Range("N31").Value = -6446
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("H74").Value = 10720.637
$ws.Range("J74").Value = 10357.375
$ws.Range("L74").Value = 10357.375
$ws.Range("N74").Value = -12229.375
$ws.Range("H77").Value = 10720.637
$ws.Range("J77").Value = 10357.375
$ws.Range("L77").Value = 31072.125
$ws.Range("N77").Value = -40432.125
$ws.Range("H81").Value = 2834.8235
$ws.Range("I81").Value = 2276
$ws.Range("J81").Value = 4176
$ws.Range("K81").Value = 4552
$ws.Range("L81").Value = 8352
$ws.Range("M81").Value = -3491
$ws.Range("N81").Value = -10474
$ws.Range("H84").Value = 2834.8235
$ws.Range("I84").Value = 2276
$ws.Range("J84").Value = 4176
$ws.Range("K84").Value = 22760
$ws.Range("L84").Value = 41760
$ws.Range("M84").Value = -17456
$ws.Range("N84").Value = -52368
